$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 25-36 (status changes / new results) ---
# Row 25
$ws.Range("C25").Value = 'Done'

# Row 26
$ws.Range("C26").Value = 'Done'
$ws.Range("E26").Value = 'X'
$ws.Range("E26").HorizontalAlignment = -4108

# Row 29
$ws.Range("C29").Value = 'Done'
$ws.Range("E29").Value = 'X'
$ws.Range("E29").HorizontalAlignment = -4108

# Row 30
$ws.Range("B30").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI None, data every 1000, restart/traj every 50,000, filename no-PI-MC'
$ws.Range("C30").Value = 'Done'
$ws.Range("E30").Value = 'X'
$ws.Range("E30").HorizontalAlignment = -4108

# Row 31
$ws.Range("B31").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 12 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MC-4'
$ws.Range("C31").Value = 'Done'
$ws.Range("E31").Value = 'X'
$ws.Range("E31").HorizontalAlignment = -4108

# Row 32
$ws.Range("B32").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 12 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MC-5'
$ws.Range("C32").Value = 'Done'
$ws.Range("E32").Value = 'X'
$ws.Range("E32").HorizontalAlignment = -4108

# Row 33
$ws.Range("B33").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 12 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MC-6'
$ws.Range("C33").Value = 'Done'
$ws.Range("E33").Value = 'X'
$ws.Range("E33").HorizontalAlignment = -4108

# Row 34
$ws.Range("B34").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.0005 ps, rCut 1.2 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MM-4'
$ws.Range("C34").Value = 'Done'
$ws.Range("E34").Value = 'X'
$ws.Range("E34").HorizontalAlignment = -4108

# Row 35
$ws.Range("B35").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.0005 ps, rCut 1.2 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MM-5'
$ws.Range("C35").Value = 'Done'
$ws.Range("E35").Value = 'X'
$ws.Range("E35").HorizontalAlignment = -4108

# Row 36
$ws.Range("B36").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.0005 ps, rCut 1.2 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MM-6'
$ws.Range("C36").Value = 'Done'
$ws.Range("E36").Value = 'X'
$ws.Range("E36").HorizontalAlignment = -4108

# --- Add new rows 37-52 ---
# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.25, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MC-7'
$ws.Range("C37").Value = 'Running'

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.25, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MC-8'
$ws.Range("C38").Value = 'Running'

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.25, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MC-9'
$ws.Range("C39").Value = 'Running'

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.00025 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MM-7'
$ws.Range("C40").Value = 'Running'

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.00025 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MM-8'
$ws.Range("C41").Value = 'Running'

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 4M steps with time step 0.00025 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MM-9'
$ws.Range("C42").Value = 'Running'

# Row 43
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 1.0, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MC-10'
$ws.Range("C43").Value = 'Running'

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 1.0, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MC-11'
$ws.Range("C44").Value = 'Pending'

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 1.0, rCut 9 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MC-12'
$ws.Range("C45").Value = 'Pending'

# Row 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 0.001 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 1 bead, data every 1000, restart/traj every 50,000, filename MM-10'
$ws.Range("C46").Value = 'Pending'

# Row 47
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 0.001 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 16 bead, data every 1000, restart/traj every 50,000, filename MM-11'
$ws.Range("C47").Value = 'Pending'

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 1M steps with time step 0.001 ps, rCut 0.9 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.997, PI 32 bead, data every 1000, restart/traj every 50,000, filename MM-12'
$ws.Range("C48").Value = 'Pending'

# Rows 49-52 (values written in the original authoring order so the
# shared-string table comes out in the same sequence as the target file)
$ws.Range("A49").Value = 48
$ws.Range("C49").Value = 'Running'
$ws.Range("D49").Value = 'NVT simulation 1 bead OpenMM'

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 9 A, NVT Andersen at 298 K, initial density 0.9979734044480487, PI 1 bead, data every 1000, restart/traj every 50,000, filename MC-NVT-1bead'
$ws.Range("C50").Value = 'Running'

$ws.Range("B49").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.0005 ps, rCut 0.9 nanometers, NVT RPMD/PILE at 298 K, initial density 0.9979734044480487, PI 1 bead, data every 1000, restart/traj every 50,000, filename MM-NVT-1bead'

$ws.Range("D50").Value = '"" in DASH'

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 'DASH simulation in dash_work/water using run_9-26-2018.sh and tip4pF_9-26-2018.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.5, rCut 15 A, NPT MonteCarlo/Andersen at 298 K and 1.0 atm, initial density 0.9979734044480487, PI 1 bead, data every 1000, restart/traj every 50,000, filename MC-rCut15'
$ws.Range("C51").Value = 'Running'
$ws.Range("D51").Value = 'Check rCut 15 in DASH'

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 'OpenMM simulation in /home/swansonk1/openmm using run_openmm.sh and pimd_modified.py, 1000 q-TIP4P/F water molecules, 2M steps with time step 0.0005 ps, rCut 1.5 nanometers, NPT RPMDMonteCarlo/PILE at 298 K and 1.01325 bar, initial density 0.9979734044480487, PI 1 bead, data every 1000, restart/traj every 50,000, filename MM-rCut15'
$ws.Range("C52").Value = 'Running'
$ws.Range("D52").Value = 'Check rCut 15 in OpenMM'

# --- Update view selection ---
$ws.Range("C53").Select()